# Fix handling of boolean fields (wordbool rather then boolean)
#
# Rows 2, 3 and 5 of the BoolCol column (F) were stored as the textual
# strings "true"/"true"/"false" instead of real Excel boolean values.
# Re-assign native PowerShell booleans so Excel stores them with
# t="b" cell typing instead of shared-string references.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = $false
$ws.Range("F3").Value = $true
$ws.Range("F5").Value = $false
